$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.546146678731126;  C = 0.546146678731126;  D = 0.3706547201795419;  E = 0.6088141918348667;  F = 0.2784751697519821;  G = 15 }
    3  = @{ B = 0.3070856474631161; C = 0.3070856474631161; D = 0.1186079915109186;  E = 0.3443951095920478;  F = 0.1617903564833405;  G = 14 }
    4  = @{ B = 0.235791455158233;  C = 0.2409394054473722; D = 0.08982400062291519; E = 0.2997065241580757;  F = 0.1925580159035662;  G = 13 }
    5  = @{ B = 0.3438942213616346; C = 0.3438942213616346; D = 0.1608659895464369;  E = 0.4010810261610949;  F = 0.2155823084169401;  G = 12 }
    6  = @{ B = 0.3554346146561582; C = 0.3554346146561582; D = 0.166620911117176;   E = 0.4081922477426244;  F = 0.2105133259524082;  G = 11 }
    7  = @{ B = 0.3119807990786818; C = 0.3119807990786818; D = 0.1298725105494941;  E = 0.3603782881216544;  F = 0.1901475788133503;  G = 10 }
    8  = @{ B = 0.342631818604201;  C = 0.342631818604201;  D = 0.1511607671542049;  E = 0.3887939906354069;  F = 0.1948967150530138;  G = 9 }
    9  = @{ B = 0.3542225829761536; C = 0.3542225829761536; D = 0.162709113486236;   E = 0.4033721773824218;  F = 0.2062882177812338;  G = 8 }
    10 = @{ B = 0.3488257250224557; C = 0.3488257250224557; D = 0.1632182784925353;  E = 0.4040028199066626;  F = 0.2201409565050135;  G = 7 }
    11 = @{ B = 0.3139331037724223; C = 0.3139331037724223; D = 0.1244752993284096;  E = 0.352810571452174;   F = 0.1763677034523832;  G = 6 }
}

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    $ws.Range("B$row").Value = $rowData.B
    $ws.Range("C$row").Value = $rowData.C
    $ws.Range("D$row").Value = $rowData.D
    $ws.Range("E$row").Value = $rowData.E
    $ws.Range("F$row").Value = $rowData.F
    $ws.Range("G$row").Value = $rowData.G
}
